# chore: update Sheets via scheduled runner
# Refresh of market-price derived figures (currentAveragePrice* / LevePrice* /
# LeveProfit* columns, H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve
# tables. A couple of rows also lose stale profit cells that the refresh no
# longer computes (cleared, not just zeroed).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1384080
$ws.Range("J17").Value = 1408547.6
$ws.Range("L17").Value = 4225642.800000001
$ws.Range("N17").Value = -4225978.800000001

$ws.Range("H38").Value = 897.9
$ws.Range("I38").Value = 331
$ws.Range("K38").Value = 993
$ws.Range("M38").Value = -621

$ws.Range("H58").Value = 1492
$ws.Range("J58").Value = 4999.5
$ws.Range("L58").Value = 14998.5
$ws.Range("N58").Value = -15298.5

$ws.Range("H107").Value = 1006.5476
$ws.Range("I107").Value = 748.5357
$ws.Range("K107").Value = 748.5357
$ws.Range("M107").Value = 1171.4643

$ws.Range("H132").Value = 55561156
$ws.Range("I132").Value = 62505924
$ws.Range("K132").Value = 187517772
$ws.Range("M132").Value = -187515242

$ws.Range("H135").Value = 2707
$ws.Range("I135").Value = 1199.7826
$ws.Range("K135").Value = 10798.0434
$ws.Range("M135").Value = -8263.0434

$ws.Range("H137").Value = 1897.0222
$ws.Range("I137").Value = 2005.3334
$ws.Range("J137").Value = 1734.5555
$ws.Range("K137").Value = 6016.0002
$ws.Range("L137").Value = 5203.666499999999
$ws.Range("M137").Value = -3466.0002
$ws.Range("N137").Value = -10303.6665

$ws.Range("H138").Value = 1536.8235
$ws.Range("I138").Value = 928.2632
$ws.Range("J138").Value = 1712.0151
$ws.Range("K138").Value = 2784.7896
$ws.Range("L138").Value = 5136.0453
$ws.Range("M138").Value = 2355.2104
$ws.Range("N138").Value = -15416.0453

$ws.Range("H141").Value = 4360.25
$ws.Range("J141").Value = 7665
$ws.Range("L141").Value = 22995
$ws.Range("N141").Value = -33355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1182.4762
$ws.Range("I2").Value = 1351.3846
$ws.Range("J2").Value = 908
$ws.Range("K2").Value = 1351.3846
$ws.Range("L2").Value = 908
$ws.Range("M2").Value = -1238.3846
$ws.Range("N2").Value = -1134

$ws.Range("H110").Value = 1443.4615
$ws.Range("I110").Value = 1286
$ws.Range("K110").Value = 1286
$ws.Range("M110").Value = 759

$ws.Range("H116").Value = 1182.4762
$ws.Range("I116").Value = 1351.3846
$ws.Range("J116").Value = 908
$ws.Range("K116").Value = 1351.3846
$ws.Range("L116").Value = 908
$ws.Range("M116").Value = 942.6153999999999
$ws.Range("N116").Value = -5496

$ws.Range("H122").Value = 3856.457
$ws.Range("I122").Value = 3421.0715
$ws.Range("J122").Value = 5598
$ws.Range("K122").Value = 10263.2145
$ws.Range("L122").Value = 16794
$ws.Range("M122").Value = -7813.2145
$ws.Range("N122").Value = -21694

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1182.4762
$ws.Range("I3").Value = 1351.3846
$ws.Range("J3").Value = 908
$ws.Range("K3").Value = 1351.3846
$ws.Range("L3").Value = 908
$ws.Range("M3").Value = -1237.3846
$ws.Range("N3").Value = -1136

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2606.1538
$ws.Range("I31").Value = 2424.4783
$ws.Range("J31").Value = 3999
$ws.Range("K31").Value = 2424.4783
$ws.Range("L31").Value = 3999
$ws.Range("M31").Value = -2129.4783
$ws.Range("N31").Value = -4589

$ws.Range("H34").Value = 2606.1538
$ws.Range("I34").Value = 2424.4783
$ws.Range("J34").Value = 3999
$ws.Range("K34").Value = 2424.4783
$ws.Range("L34").Value = 3999
$ws.Range("M34").Value = -2222.4783
$ws.Range("N34").Value = -4403

$ws.Range("H58").Value = 3555.348
$ws.Range("I58").Value = 2897.25
$ws.Range("J58").Value = 3906.3333
$ws.Range("K58").Value = 2897.25
$ws.Range("L58").Value = 3906.3333
$ws.Range("M58").Value = -2694.25
$ws.Range("N58").Value = -4312.3333

$ws.Range("H76").Value = 10009.75
$ws.Range("I76").Value = 10009.75
$ws.Range("K76").Value = 10009.75
$ws.Range("M76").Value = -9694.75

$ws.Range("H79").Value = 10009.75
$ws.Range("I79").Value = 10009.75
$ws.Range("K79").Value = 10009.75
$ws.Range("M79").Value = -8917.75

$ws.Range("H99").Value = 3332.2334
$ws.Range("I99").Value = 2713.7827
$ws.Range("J99").Value = 5364.2856
$ws.Range("K99").Value = 2713.7827
$ws.Range("L99").Value = 5364.2856
$ws.Range("M99").Value = -1215.7827
$ws.Range("N99").Value = -8360.285599999999

$ws.Range("H107").Value = 541.7143
$ws.Range("I107").Value = 324.25
$ws.Range("K107").Value = 324.25
$ws.Range("M107").Value = 1595.75

$ws.Range("H126").Value = 3332.2334
$ws.Range("I126").Value = 2713.7827
$ws.Range("J126").Value = 5364.2856
$ws.Range("K126").Value = 8141.348100000001
$ws.Range("L126").Value = 16092.8568
$ws.Range("M126").Value = -5671.348100000001
$ws.Range("N126").Value = -21032.8568

$ws.Range("H134").Value = 4634.9565
$ws.Range("I134").Value = 2496.5386
$ws.Range("K134").Value = 7489.6158
$ws.Range("M134").Value = -4954.6158

$ws.Range("H136").Value = 3555.348
$ws.Range("I136").Value = 2897.25
$ws.Range("J136").Value = 3906.3333
$ws.Range("K136").Value = 8691.75
$ws.Range("L136").Value = 11718.9999
$ws.Range("M136").Value = -6141.75
$ws.Range("N136").Value = -16818.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 136.53334
$ws.Range("I2").Value = 149.55
$ws.Range("J2").Value = 110.5
$ws.Range("K2").Value = 897.3000000000001
$ws.Range("L2").Value = 663
$ws.Range("M2").Value = -784.3000000000001
$ws.Range("N2").Value = -889

$ws.Range("H38").Value = 4068.75
$ws.Range("I38").Value = 417.14285
$ws.Range("J38").Value = 6908.8887
$ws.Range("K38").Value = 1251.42855
$ws.Range("L38").Value = 20726.6661
$ws.Range("M38").Value = -904.4285500000001
$ws.Range("N38").Value = -21420.6661

$ws.Range("H92").Value = 861.9091
$ws.Range("J92").Value = 918.44446
$ws.Range("L92").Value = 2755.33338
$ws.Range("N92").Value = -5251.33338

$ws.Range("H113").Value = 2001.3846
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2001.3846
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6004.1538
$ws.Range("N113").Value = -10344.1538
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 129580.625
$ws.Range("I126").Value = 501999.5
$ws.Range("J126").Value = 5441
$ws.Range("K126").Value = 1505998.5
$ws.Range("L126").Value = 16323
$ws.Range("M126").Value = -1503528.5
$ws.Range("N126").Value = -21263

$ws.Range("H136").Value = 41468.434
$ws.Range("J136").Value = 41468.434
$ws.Range("L136").Value = 124405.302
$ws.Range("N136").Value = -129505.302

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3621.6086
$ws.Range("I61").Value = 3745.0908
$ws.Range("J61").Value = 905
$ws.Range("K61").Value = 3745.0908
$ws.Range("L61").Value = 905
$ws.Range("M61").Value = -3543.0908
$ws.Range("N61").Value = -1309

$ws.Range("H113").Value = 3621.6086
$ws.Range("I113").Value = 3745.0908
$ws.Range("K113").Value = 3745.0908
$ws.Range("M113").Value = -1575.0908
$ws.Range("N113").Value = -5245

$ws.Range("H132").Value = 3065.4243
$ws.Range("J132").Value = 3095.9546
$ws.Range("L132").Value = 9287.863799999999
$ws.Range("N132").Value = -14347.8638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1046.6111
$ws.Range("I107").Value = 883.8461
$ws.Range("K107").Value = 2651.5383
$ws.Range("M107").Value = -731.5383000000002

$ws.Range("H113").Value = 5560878.5
$ws.Range("I113").Value = 7582592.5
$ws.Range("J113").Value = 1165.5
$ws.Range("K113").Value = 22747777.5
$ws.Range("L113").Value = 3496.5
$ws.Range("M113").Value = -22745607.5
$ws.Range("N113").Value = -7836.5

$ws.Range("H122").Value = 2115.1365
$ws.Range("I122").Value = 2058.4375
$ws.Range("K122").Value = 6175.3125
$ws.Range("M122").Value = -3725.3125

$ws.Range("H126").Value = 1965.76
$ws.Range("I126").Value = 1479.9166
$ws.Range("K126").Value = 4439.7498
$ws.Range("M126").Value = -1969.7498

$ws.Range("H132").Value = 1211.3948
$ws.Range("I132").Value = 1104.375
$ws.Range("K132").Value = 3313.125
$ws.Range("M132").Value = -783.125
